# Apply scraped cryptocurrency price/volume updates from the Dec 8, 2023
# GitHub Actions run. All data cells in this sheet are stored as plain text
# (inlineStr) in the source workbook. Values that look like plain numbers
# (e.g. "3.80", "172.20") are written with a leading apostrophe so Excel
# keeps them as text instead of silently re-interpreting them as numbers
# (which would drop trailing zeros / change the cell type).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.667.94"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3
$ws.Range("D3").Value = "2.351.63"
$ws.Range("E3").Value = "  +4.58%  "

# Row 4
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").Value = "'234.91"
$ws.Range("E5").Value = "  +1.93%  "

# Row 6
$ws.Range("E6").Value = "  +2.71%  "

# Row 7
$ws.Range("D7").Value = "'73.21"
$ws.Range("E7").Value = "  +14.30%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  +19.57%  "

# Row 10
$ws.Range("D10").Value = "'0.0982"
$ws.Range("E10").Value = "  +3.33%  "

# Row 11
$ws.Range("D11").Value = "'27.36"
$ws.Range("E11").Value = "  +3.21%  "

# Row 12
$ws.Range("D12").Value = "'0.106"
$ws.Range("E12").Value = "  +2.45%  "

# Row 13
$ws.Range("D13").Value = "2.701.69"
$ws.Range("E13").Value = "  +4.61%  "

# Row 14
$ws.Range("D14").Value = "'16.89"
$ws.Range("E14").Value = "  +13.71%  "

# Row 15
$ws.Range("D15").Value = "'6.64"
$ws.Range("E15").Value = "  +10.40%  "

# Row 16
$ws.Range("D16").Value = "'0.878"
$ws.Range("E16").Value = "  +7.29%  "

# Row 17
$ws.Range("D17").Value = "2.356.36"
$ws.Range("E17").Value = "  +4.80%  "

# Row 18
$ws.Range("D18").Value = "43.541.27"
$ws.Range("E18").Value = "  +0.76%  "

# Row 19
$ws.Range("D19").Value = "'0.0000101"
$ws.Range("E19").Value = "  +4.46%  "

# Row 20
$ws.Range("D20").Value = "'75.78"
$ws.Range("E20").Value = "  +3.93%  "

# Row 21
$ws.Range("D21").Value = "'6.41"
$ws.Range("E21").Value = "  +5.94%  "

# Row 22
$ws.Range("D22").Value = "'250.55"
$ws.Range("E22").Value = "  +1.60%  "

# Row 23
$ws.Range("D23").Value = "'3.80"
$ws.Range("E23").Value = "  -2.48%  "

# Row 25
$ws.Range("E25").Value = "  +2.00%  "

# Row 26
$ws.Range("D26").Value = "'10.19"
$ws.Range("E26").Value = "  +5.11%  "

# Row 27
$ws.Range("D27").Value = "'2.25"
$ws.Range("E27").Value = "  -1.88%  "

# Row 28
$ws.Range("D28").Value = "'22.41"
$ws.Range("E28").Value = "  +4.01%  "

# Row 29
$ws.Range("D29").Value = "'172.20"
$ws.Range("E29").Value = "  -0.85%  "

# Row 30
$ws.Range("E30").Value = "  +7.87%  "

# Row 31
$ws.Range("E31").Value = "  +2.85%  "

# Row 32
$ws.Range("E32").Value = "  +4.35%  "

# Row 33
$ws.Range("D33").Value = "'5.08"
$ws.Range("E33").Value = "  +3.43%  "

# Row 34
$ws.Range("E34").Value = "  +3.30%  "

# Row 35
$ws.Range("D35").Value = "'5.09"
$ws.Range("E35").Value = "  +3.50%  "

# Row 36
$ws.Range("D36").Value = "'3.74"
$ws.Range("E36").Value = "  +3.16%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.43"
$ws.Range("E37").Value = "  +7.58%  "

# Row 38
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "'6.44"
$ws.Range("E38").Value = "  +1.63%  "

# Row 39
$ws.Range("D39").Value = "'0.0265"
$ws.Range("E39").Value = "  +6.64%  "

# Row 40
$ws.Range("D40").Value = "'19.47"
$ws.Range("E40").Value = "  +13.45%  "

# Row 41
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("E42").Value = "  +0.61%  "

# Row 43
$ws.Range("E43").Value = "  +9.12%  "

# Row 44
$ws.Range("D44").Value = "'98.98"
$ws.Range("E44").Value = "  +2.62%  "

# Row 45
$ws.Range("D45").Value = "'1.22"
$ws.Range("E45").Value = "  +3.48%  "

# Row 46
$ws.Range("D46").Value = "'0.0960"
$ws.Range("E46").Value = "  +2.80%  "

# Row 47
$ws.Range("D47").Value = "'4.42"
$ws.Range("E47").Value = "  -1.48%  "

# Row 48
$ws.Range("D48").Value = "'0.181"
$ws.Range("E48").Value = "  +13.51%  "

# Row 49
$ws.Range("D49").Value = "1.438.85"
$ws.Range("E49").Value = "  +0.65%  "

# Row 50
$ws.Range("D50").Value = "2.579.12"
$ws.Range("E50").Value = "  +4.73%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'2.29"
$ws.Range("E51").Value = "  +1.50%  "
